$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()
$win = $excel.ActiveWindow
Write-Host "trying LargeScroll"
$r = $win.LargeScroll(1,0,0,0)
Write-Host "result: $r"
